$wb = $excel.ActiveWorkbook

# The "Nějaký import" worksheet held a couple of DTO-resolution related rows
# that moved elsewhere - clear the now-unused cells and remove the stray
# "Bad one" row.
$ws = $wb.Worksheets.Item("Nějaký import")

# Row 2 (B2:C2) used to carry "Another Value 001" / "Va-Va-Value 001" - clear it.
$ws.Range("B2:C2").ClearContents()

# The very last row only held a leftover "Bad one" marker - remove it entirely.
$ws.Rows.Item(22).Delete()

# Restore/point the selection on the "tabs" sheet (this also was nudged in the
# original edit) before switching the active sheet, since selecting a range
# on a worksheet also activates that worksheet.
$tabsWs = $wb.Worksheets.Item("tabs")
$tabsWs.Range("B3").Select()

# Select the new working cell on "Nějaký import" and make it the active sheet.
$ws.Range("D6").Select()

# Scroll the sheet tab strip so it is brought into view.
$excel.ActiveWindow.ScrollWorkbookTabs(1)
